$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell (row 464, column A) carries the date-style formatting (border, bold, centered, custom date number format).
# We replicate that formatting onto each newly appended date cell via a format-only paste so the workbook reuses the
# existing style index instead of registering a brand new one.

$newRows = @(
    @(465, 44539, 71, 186, 257.0161256891763),
    @(466, 44540, 49, 210, 290.1794967458442),
    @(467, 44541, 7, 192, 265.3069684533432),
    @(468, 44542, 34, 202, 279.1250397269549),
    @(469, 44543, 34, 204, 281.8886539816772),
    @(470, 44544, 37, 234, 323.3428678025121),
    @(471, 44545, 0, 232, 320.5792535477898),
    @(472, 44546, 25, 186, 257.0161256891763),
    @(473, 44547, 58, 195, 269.4523898354267),
    @(474, 44548, 22, 210, 290.1794967458442),
    @(475, 44550, 42, 218, 301.2339537647335),
    @(476, 44551, 58, 242, 334.3973248214014),
    @(477, 44552, 7, 212, 292.9431110005665),
    @(478, 44553, 30, 242, 334.3973248214014),
    @(479, 44554, 16, 233, 321.961060675151),
    @(480, 44555, 63, 238, 328.8700963119568),
    @(481, 44556, 19, 235, 324.7246749298733),
    @(482, 44557, 72, 265, 366.1788887507082),
    @(483, 44558, 68, 275, 379.9969600243198),
    @(484, 44559, 81, 349, 482.2506874490459),
    @(485, 44560, 20, 339, 468.4326161754342),
    @(486, 44561, 57, 380, 525.086708397242),
    @(487, 44562, 31, 348, 480.8688803216847),
    @(488, 44563, 166, 495, 683.9945280437756),
    @(489, 44564, 174, 597, 824.9388550346143),
    @(490, 44565, 53, 582, 804.2117481241968),
    @(491, 44566, 66, 567, 783.4846412137794)
)

foreach ($entry in $newRows) {
    $rowIndex = $entry[0]
    $dateSerial = $entry[1]
    $nuoviPos = $entry[2]
    $sommaMobile = $entry[3]
    $incidenza = $entry[4]

    $ws.Cells.Item(464, 1).Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value2 = $dateSerial
    $ws.Cells.Item($rowIndex, 2).Value2 = $nuoviPos
    $ws.Cells.Item($rowIndex, 3).Value2 = $sommaMobile
    $ws.Cells.Item($rowIndex, 4).Value2 = $incidenza
}

$excel.CutCopyMode = $false

Write-Host "Appended $($newRows.Count) rows; new used range: $($ws.UsedRange.Address())"
